$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("By course")

# Intro CS courses (COMP_SCI 101, 110, 111, 211, 230, 295) get the new
# "ComputerScience" topic added in column F (Course Topic Area).
$rows = @(2, 3, 4, 5, 6, 7, 8, 10, 11, 12, 25, 27)

foreach ($r in $rows) {
    $ws.Range("F$r").Value = "ComputerScience"
}

$ws.Range("F25").Select()
